# Update the word-level bounding box table: split the previous multi-word
# cells into one row per word, with updated x/y/width/height values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data: column A text, B = x, C = y, D = width, E = height
$data = @(
    @("Greek ",        152, 644,    58,  23),
    @("mythology ",     210, 644,    93,  23),
    @("gift ",          870, 926.8,  31,  23),
    @("of ",            901, 926.8,  22,  23),
    @("prophecy. ",     152, 951.8,  89,  23),
    @("Trojan ",        356, 1001.8, 59,  23),
    @("Horse ",         415, 1001.8, 57,  23),
    @("trick, ",        472, 1001.8, 46,  23),
    @("Agamemnon's ",   226, 1026.8, 128, 23),
    @("Bronze ",        349, 1193.2, 67,  23),
    @("Age. ",          416, 1193.2, 45,  23)
)

# Fix the apostrophe to the typographic right single quotation mark (’)
$data[8][0] = "Agamemnon" + [char]0x2019 + "s "

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
